$wb = $excel.ActiveWorkbook

# --- Logs sheet: append a new row (row 3) with the latest log entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A3").Value = "Demo inplannen"
$ws.Range("B3").Value = "planning@testbedrijf123.nl"
$ws.Range("C3").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D3").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E3").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F3").Value = "2025-08-14 19:05:16"
$ws.Range("G3").Value = "Nee"
$ws.Range("H3").Value = "Ja"
$ws.Range("I3").Value = "Nee"
$ws.Range("J3").Value = "Nee"

# Extend the per-column conditional formatting so it also covers row 3
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $fcs = $ws.Range($col + "2").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range($col + "2:" + $col + "3"))
    }
}

# --- Dashboard sheet: bump the count for the existing category ---
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B2").Value = 2
